$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.119.02"
$ws.Range("E2").Value = "  +2.29%  "

$ws.Range("D3").Value = "3.377.67"
$ws.Range("E3").Value = "  +7.58%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'260.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.56%  "

$ws.Range("D6").Value = "'629.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("E7").Value = "  +23.10%  "

$ws.Range("D8").Value = "'0.394"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.26%  "

$ws.Range("D10").Value = "'0.864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.34%  "

$ws.Range("D11").Value = "3.375.96"
$ws.Range("E11").Value = "  +7.63%  "

$ws.Range("D12").Value = "'0.199"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "98.826.12"
$ws.Range("E13").Value = "  +2.21%  "

$ws.Range("D14").Value = "'36.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.32%  "

$ws.Range("D15").Value = "'0.0000249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.11%  "

$ws.Range("D16").Value = "4.006.45"
$ws.Range("E16").Value = "  +7.79%  "

$ws.Range("D17").Value = "'5.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").Value = "3.366.22"
$ws.Range("E18").Value = "  +7.34%  "

$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").Value = "'15.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.71%  "

$ws.Range("D21").Value = "'494.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.31%  "

$ws.Range("D22").Value = "'6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.44%  "

$ws.Range("D23").Value = "'0.0000211"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.57%  "

$ws.Range("D24").Value = "'9.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.09%  "

$ws.Range("D25").Value = "'5.65"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'88.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "'11.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("D28").Value = "3.555.61"
$ws.Range("E28").Value = "  +7.69%  "

$ws.Range("D29").Value = "'0.281"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.64%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").Value = "'0.193"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.52%  "

$ws.Range("D32").Value = "'0.131"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.55%  "

$ws.Range("D33").Value = "'0.996"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").Value = "'9.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.05%  "

$ws.Range("D35").Value = "'28.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.16%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.151"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "'1.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.15%  "

$ws.Range("B39").Value = "MantraDAO"
$ws.Range("C39").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D39").Value = "'4.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.29%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'500.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.37%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "'0.461"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.96%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'24.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.97%  "

$ws.Range("D43").Value = "'1.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("D44").Value = "'3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.59%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.786"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.33%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "'161.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("D49").Value = "'0.842"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +13.45%  "

$ws.Range("E50").Value = "  +2.71%  "

$ws.Range("D51").Value = "'46.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.09%  "
